$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 34.78985266666667
$ws.Range("H2").Value = 104.369558
$ws.Range("I2").Value = 0.09964961663893999
$ws.Range("J2").Value = 0.09964961663893998
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 3.827857909208
$ws.Range("R2").Value = 34.45072118287199
$ws.Range("S2").Value = 0.09964961663893999
$ws.Range("T2").Value = 0.09964961663893998

# Row 3
$ws.Range("I3").Value = 0.1313371392780071
$ws.Range("J3").Value = 0.1313371392780071
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("Q3").Value = 5.045076181072
$ws.Range("R3").Value = 45.405685629648
$ws.Range("S3").Value = 0.1313371392780071
$ws.Range("T3").Value = 0.1313371392780071

# Row 4
$ws.Range("G4").Value = 32.884922
$ws.Range("H4").Value = 98.654766
$ws.Range("I4").Value = 0.09419326669472271
$ws.Range("J4").Value = 0.09419326669472271
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("Q4").Value = 3.618262197816
$ws.Range("R4").Value = 32.564359780344
$ws.Range("S4").Value = 0.09419326669472271
$ws.Range("T4").Value = 0.09419326669472271

# Row 5
$ws.Range("G5").Value = 7.302655333333334
$ws.Range("H5").Value = 21.907966
$ws.Range("I5").Value = 0.02091721432066362
$ws.Range("J5").Value = 0.02091721432066361
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.110028
$ws.Range("N5").Value = 0.330084
$ws.Range("Q5").Value = 0.8034965610160001
$ws.Range("R5").Value = 7.231469049144001
$ws.Range("S5").Value = 0.02091721432066362
$ws.Range("T5").Value = 0.02091721432066361

# Row 6
$ws.Range("G6").Value = 28.09959066666667
$ws.Range("H6").Value = 84.298772
$ws.Range("I6").Value = 0.08048649887866162
$ws.Range("J6").Value = 0.08048649887866162
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.110028
$ws.Range("N6").Value = 0.330084
$ws.Range("Q6").Value = 3.091741761872
$ws.Range("R6").Value = 27.825675856848
$ws.Range("S6").Value = 0.08048649887866162
$ws.Range("T6").Value = 0.08048649887866162

# Row 7
$ws.Range("G7").Value = 200.1921133333333
$ws.Range("H7").Value = 600.57634
$ws.Range("I7").Value = 0.5734162641890049
$ws.Range("J7").Value = 0.5734162641890049
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.110028
$ws.Range("N7").Value = 0.330084
$ws.Range("Q7").Value = 22.02673784584
$ws.Range("R7").Value = 198.24064061256
$ws.Range("S7").Value = 0.5734162641890049
$ws.Range("T7").Value = 0.5734162641890049
